$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) PairCorrInput sheet ("PairCorrInput") - update asset/api parameters
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("PairCorrInput")

# Row 5: Mode -> api (was disk); clear the now-unused period value in F5
# (set first so the shared-string table picks up 'api' ahead of the new
# asset/api tickers below, matching the order they were authored in)
$ws1.Range("B5").Value = "api"
$ws1.Range("F5").ClearContents()

# Row 2: asset1 / api1 / asset2 / api2 / NumDays / CC(1) -- switch sample pair
# from BTC/coingecko vs XAUUSD,FX_IDC/tv to BTCUSD,INDEX/tv vs NQ1!,CME_MINI/tv
$ws1.Range("C2").Value = "NQ1!,CME_MINI"
$ws1.Range("A2").Value = "BTCUSD,INDEX"
$ws1.Range("B2").Value = "tv"
$ws1.Range("D2").Value = "tv"
$ws1.Range("E2").Value = 250
$ws1.Range("F2").Value = 30
$ws1.Range("F2").ClearFormats()

# Row 3/4: CC period values (days)
$ws1.Range("F3").Value = 90
$ws1.Range("F4").Value = 180

# New formula cell: 10 years in days
$ws1.Range("I24").Formula = "=10*365"

# ---------------------------------------------------------------------------
# 2) Information sheet - update explanatory text for the price/yoy parameter
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Information")
$ws2.Range("C24").Value = "price'or 'yoy'to display assets as price or YoY % change (first derivative of price wrt time, using rolling period of 1 year)"

# ---------------------------------------------------------------------------
# 3) Add a new "Sheet1" worksheet at the end with a column of monthly dates
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add()
$newSheetName = $newSheet.Name
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Move($null, $lastSheet)
$ws3 = $wb.Worksheets.Item($newSheetName)

$ws3.Range("E4:E22").NumberFormat = "yyyy-mm-dd;@"
$dates = @(44927,44958,44986,45017,45047,45078,45108,45139,45170,45200,45231,45261,45292)
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = 4 + $i
    $ws3.Cells.Item($row, 5).Value = $dates[$i]
}
$ws3.Range("D3").Select()

# ---------------------------------------------------------------------------
# 4) Restore tab/selection state: PairCorrInput active with E10 selected,
#    Information no longer the active tab, selection moved to D29.
# ---------------------------------------------------------------------------
$ws2.Range("D29").Select()
$ws1.Activate()
$ws1.Range("E10").Select()
